$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the benchmark data for row 14 (Ungoogled Chromium)
$ws.Range("B14").Value = 210
$ws.Range("C14").Value = 1379
$ws.Range("D14").Value = 198
$ws.Range("E14").Value = 1760

# Formula for the average column, consistent with the rest of the table
$ws.Range("F14").Formula = "=AVERAGE(B14:E14)"

# Update the active selection to F14
$ws.Range("F14").Select()
